$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 39-58 ---
$ws = $wb.Worksheets.Item("PIR")
$newRange = $ws.Range("A39:F58")
$newRange.NumberFormat = "@"

$ws.Cells.Item(39, 1).Value = '2026-02-01'
$ws.Cells.Item(39, 2).Value = '13:48:25'
$ws.Cells.Item(39, 3).Value = '13:00'
$ws.Cells.Item(39, 4).Value = 'Bathroom'
$ws.Cells.Item(39, 5).Value = 'No Motion'
$ws.Cells.Item(39, 6).Value = 'Inactive'

$ws.Cells.Item(40, 1).Value = '2026-02-01'
$ws.Cells.Item(40, 2).Value = '13:48:26'
$ws.Cells.Item(40, 3).Value = '13:00'
$ws.Cells.Item(40, 4).Value = 'Bathroom'
$ws.Cells.Item(40, 5).Value = 'No Motion'
$ws.Cells.Item(40, 6).Value = 'Inactive'

$ws.Cells.Item(41, 1).Value = '2026-02-01'
$ws.Cells.Item(41, 2).Value = '13:48:28'
$ws.Cells.Item(41, 3).Value = '13:00'
$ws.Cells.Item(41, 4).Value = 'Bathroom'
$ws.Cells.Item(41, 5).Value = 'Motion Detected'
$ws.Cells.Item(41, 6).Value = 'Active'

$ws.Cells.Item(42, 1).Value = '2026-02-01'
$ws.Cells.Item(42, 2).Value = '13:48:28'
$ws.Cells.Item(42, 3).Value = '13:00'
$ws.Cells.Item(42, 4).Value = 'Bathroom'
$ws.Cells.Item(42, 5).Value = 'Motion Detected'
$ws.Cells.Item(42, 6).Value = 'Active'

$ws.Cells.Item(43, 1).Value = '2026-02-01'
$ws.Cells.Item(43, 2).Value = '13:48:28'
$ws.Cells.Item(43, 3).Value = '13:00'
$ws.Cells.Item(43, 4).Value = 'Bathroom'
$ws.Cells.Item(43, 5).Value = 'No Motion'
$ws.Cells.Item(43, 6).Value = 'Inactive'

$ws.Cells.Item(44, 1).Value = '2026-02-01'
$ws.Cells.Item(44, 2).Value = '13:48:29'
$ws.Cells.Item(44, 3).Value = '13:00'
$ws.Cells.Item(44, 4).Value = 'Bathroom'
$ws.Cells.Item(44, 5).Value = 'No Motion'
$ws.Cells.Item(44, 6).Value = 'Inactive'

$ws.Cells.Item(45, 1).Value = '2026-02-01'
$ws.Cells.Item(45, 2).Value = '13:48:29'
$ws.Cells.Item(45, 3).Value = '13:00'
$ws.Cells.Item(45, 4).Value = 'Bathroom'
$ws.Cells.Item(45, 5).Value = 'Motion Detected'
$ws.Cells.Item(45, 6).Value = 'Active'

$ws.Cells.Item(46, 1).Value = '2026-02-01'
$ws.Cells.Item(46, 2).Value = '13:48:29'
$ws.Cells.Item(46, 3).Value = '13:00'
$ws.Cells.Item(46, 4).Value = 'Bathroom'
$ws.Cells.Item(46, 5).Value = 'Motion Detected'
$ws.Cells.Item(46, 6).Value = 'Active'

$ws.Cells.Item(47, 1).Value = '2026-02-01'
$ws.Cells.Item(47, 2).Value = '13:48:30'
$ws.Cells.Item(47, 3).Value = '13:00'
$ws.Cells.Item(47, 4).Value = 'Bathroom'
$ws.Cells.Item(47, 5).Value = 'No Motion'
$ws.Cells.Item(47, 6).Value = 'Inactive'

$ws.Cells.Item(48, 1).Value = '2026-02-01'
$ws.Cells.Item(48, 2).Value = '13:48:30'
$ws.Cells.Item(48, 3).Value = '13:00'
$ws.Cells.Item(48, 4).Value = 'Bathroom'
$ws.Cells.Item(48, 5).Value = 'No Motion'
$ws.Cells.Item(48, 6).Value = 'Inactive'

$ws.Cells.Item(49, 1).Value = '2026-02-01'
$ws.Cells.Item(49, 2).Value = '13:48:30'
$ws.Cells.Item(49, 3).Value = '13:00'
$ws.Cells.Item(49, 4).Value = 'Bathroom'
$ws.Cells.Item(49, 5).Value = 'Motion Detected'
$ws.Cells.Item(49, 6).Value = 'Active'

$ws.Cells.Item(50, 1).Value = '2026-02-01'
$ws.Cells.Item(50, 2).Value = '13:48:30'
$ws.Cells.Item(50, 3).Value = '13:00'
$ws.Cells.Item(50, 4).Value = 'Bathroom'
$ws.Cells.Item(50, 5).Value = 'Motion Detected'
$ws.Cells.Item(50, 6).Value = 'Active'

$ws.Cells.Item(51, 1).Value = '2026-02-01'
$ws.Cells.Item(51, 2).Value = '13:48:31'
$ws.Cells.Item(51, 3).Value = '13:00'
$ws.Cells.Item(51, 4).Value = 'Bathroom'
$ws.Cells.Item(51, 5).Value = 'No Motion'
$ws.Cells.Item(51, 6).Value = 'Inactive'

$ws.Cells.Item(52, 1).Value = '2026-02-01'
$ws.Cells.Item(52, 2).Value = '13:48:32'
$ws.Cells.Item(52, 3).Value = '13:00'
$ws.Cells.Item(52, 4).Value = 'Bathroom'
$ws.Cells.Item(52, 5).Value = 'No Motion'
$ws.Cells.Item(52, 6).Value = 'Inactive'

$ws.Cells.Item(53, 1).Value = '2026-02-01'
$ws.Cells.Item(53, 2).Value = '13:48:32'
$ws.Cells.Item(53, 3).Value = '13:00'
$ws.Cells.Item(53, 4).Value = 'Bathroom'
$ws.Cells.Item(53, 5).Value = 'No Motion'
$ws.Cells.Item(53, 6).Value = 'Inactive'

$ws.Cells.Item(54, 1).Value = '2026-02-01'
$ws.Cells.Item(54, 2).Value = '13:48:32'
$ws.Cells.Item(54, 3).Value = '13:00'
$ws.Cells.Item(54, 4).Value = 'Bathroom'
$ws.Cells.Item(54, 5).Value = 'No Motion'
$ws.Cells.Item(54, 6).Value = 'Inactive'

$ws.Cells.Item(55, 1).Value = '2026-02-01'
$ws.Cells.Item(55, 2).Value = '13:48:35'
$ws.Cells.Item(55, 3).Value = '13:00'
$ws.Cells.Item(55, 4).Value = 'Bathroom'
$ws.Cells.Item(55, 5).Value = 'No Motion'
$ws.Cells.Item(55, 6).Value = 'Inactive'

$ws.Cells.Item(56, 1).Value = '2026-02-01'
$ws.Cells.Item(56, 2).Value = '13:48:36'
$ws.Cells.Item(56, 3).Value = '13:00'
$ws.Cells.Item(56, 4).Value = 'Bathroom'
$ws.Cells.Item(56, 5).Value = 'No Motion'
$ws.Cells.Item(56, 6).Value = 'Inactive'

$ws.Cells.Item(57, 1).Value = '2026-02-01'
$ws.Cells.Item(57, 2).Value = '13:48:40'
$ws.Cells.Item(57, 3).Value = '13:00'
$ws.Cells.Item(57, 4).Value = 'Bathroom'
$ws.Cells.Item(57, 5).Value = 'No Motion'
$ws.Cells.Item(57, 6).Value = 'Inactive'

$ws.Cells.Item(58, 1).Value = '2026-02-01'
$ws.Cells.Item(58, 2).Value = '13:48:41'
$ws.Cells.Item(58, 3).Value = '13:00'
$ws.Cells.Item(58, 4).Value = 'Bathroom'
$ws.Cells.Item(58, 5).Value = 'No Motion'
$ws.Cells.Item(58, 6).Value = 'Inactive'

$newRange.ClearFormats()

# --- Humidity sheet: append rows 18-33 ---
$ws = $wb.Worksheets.Item("Humidity")
$newRange = $ws.Range("A18:F33")
$newRange.NumberFormat = "@"

$ws.Cells.Item(18, 1).Value = '2026-02-01'
$ws.Cells.Item(18, 2).Value = '13:48:22'
$ws.Cells.Item(18, 3).Value = '13:00'
$ws.Cells.Item(18, 4).Value = 'Bathroom'
$ws.Cells.Item(18, 5).Value = '99.9%'
$ws.Cells.Item(18, 6).Value = 'Active'

$ws.Cells.Item(19, 1).Value = '2026-02-01'
$ws.Cells.Item(19, 2).Value = '13:48:23'
$ws.Cells.Item(19, 3).Value = '13:00'
$ws.Cells.Item(19, 4).Value = 'Bathroom'
$ws.Cells.Item(19, 5).Value = '99.9%'
$ws.Cells.Item(19, 6).Value = 'Active'

$ws.Cells.Item(20, 1).Value = '2026-02-01'
$ws.Cells.Item(20, 2).Value = '13:48:27'
$ws.Cells.Item(20, 3).Value = '13:00'
$ws.Cells.Item(20, 4).Value = 'Bathroom'
$ws.Cells.Item(20, 5).Value = '99.8%'
$ws.Cells.Item(20, 6).Value = 'Active'

$ws.Cells.Item(21, 1).Value = '2026-02-01'
$ws.Cells.Item(21, 2).Value = '13:48:28'
$ws.Cells.Item(21, 3).Value = '13:00'
$ws.Cells.Item(21, 4).Value = 'Bathroom'
$ws.Cells.Item(21, 5).Value = '80.8%'
$ws.Cells.Item(21, 6).Value = 'Active'

$ws.Cells.Item(22, 1).Value = '2026-02-01'
$ws.Cells.Item(22, 2).Value = '13:48:28'
$ws.Cells.Item(22, 3).Value = '13:00'
$ws.Cells.Item(22, 4).Value = 'Bathroom'
$ws.Cells.Item(22, 5).Value = '80.8%'
$ws.Cells.Item(22, 6).Value = 'Active'

$ws.Cells.Item(23, 1).Value = '2026-02-01'
$ws.Cells.Item(23, 2).Value = '13:48:29'
$ws.Cells.Item(23, 3).Value = '13:00'
$ws.Cells.Item(23, 4).Value = 'Bathroom'
$ws.Cells.Item(23, 5).Value = '80.6%'
$ws.Cells.Item(23, 6).Value = 'Active'

$ws.Cells.Item(24, 1).Value = '2026-02-01'
$ws.Cells.Item(24, 2).Value = '13:48:29'
$ws.Cells.Item(24, 3).Value = '13:00'
$ws.Cells.Item(24, 4).Value = 'Bathroom'
$ws.Cells.Item(24, 5).Value = '80.8%'
$ws.Cells.Item(24, 6).Value = 'Active'

$ws.Cells.Item(25, 1).Value = '2026-02-01'
$ws.Cells.Item(25, 2).Value = '13:48:30'
$ws.Cells.Item(25, 3).Value = '13:00'
$ws.Cells.Item(25, 4).Value = 'Bathroom'
$ws.Cells.Item(25, 5).Value = '81.6%'
$ws.Cells.Item(25, 6).Value = 'Active'

$ws.Cells.Item(26, 1).Value = '2026-02-01'
$ws.Cells.Item(26, 2).Value = '13:48:31'
$ws.Cells.Item(26, 3).Value = '13:00'
$ws.Cells.Item(26, 4).Value = 'Bathroom'
$ws.Cells.Item(26, 5).Value = '99.9%'
$ws.Cells.Item(26, 6).Value = 'Active'

$ws.Cells.Item(27, 1).Value = '2026-02-01'
$ws.Cells.Item(27, 2).Value = '13:48:31'
$ws.Cells.Item(27, 3).Value = '13:00'
$ws.Cells.Item(27, 4).Value = 'Bathroom'
$ws.Cells.Item(27, 5).Value = '99.9%'
$ws.Cells.Item(27, 6).Value = 'Active'

$ws.Cells.Item(28, 1).Value = '2026-02-01'
$ws.Cells.Item(28, 2).Value = '13:48:31'
$ws.Cells.Item(28, 3).Value = '13:00'
$ws.Cells.Item(28, 4).Value = 'Bathroom'
$ws.Cells.Item(28, 5).Value = '99.9%'
$ws.Cells.Item(28, 6).Value = 'Active'

$ws.Cells.Item(29, 1).Value = '2026-02-01'
$ws.Cells.Item(29, 2).Value = '13:48:32'
$ws.Cells.Item(29, 3).Value = '13:00'
$ws.Cells.Item(29, 4).Value = 'Bathroom'
$ws.Cells.Item(29, 5).Value = '99.8%'
$ws.Cells.Item(29, 6).Value = 'Active'

$ws.Cells.Item(30, 1).Value = '2026-02-01'
$ws.Cells.Item(30, 2).Value = '13:48:32'
$ws.Cells.Item(30, 3).Value = '13:00'
$ws.Cells.Item(30, 4).Value = 'Bathroom'
$ws.Cells.Item(30, 5).Value = '92.4%'
$ws.Cells.Item(30, 6).Value = 'Active'

$ws.Cells.Item(31, 1).Value = '2026-02-01'
$ws.Cells.Item(31, 2).Value = '13:48:33'
$ws.Cells.Item(31, 3).Value = '13:00'
$ws.Cells.Item(31, 4).Value = 'Bathroom'
$ws.Cells.Item(31, 5).Value = '90.0%'
$ws.Cells.Item(31, 6).Value = 'Active'

$ws.Cells.Item(32, 1).Value = '2026-02-01'
$ws.Cells.Item(32, 2).Value = '13:48:37'
$ws.Cells.Item(32, 3).Value = '13:00'
$ws.Cells.Item(32, 4).Value = 'Bathroom'
$ws.Cells.Item(32, 5).Value = '89.1%'
$ws.Cells.Item(32, 6).Value = 'Active'

$ws.Cells.Item(33, 1).Value = '2026-02-01'
$ws.Cells.Item(33, 2).Value = '13:48:38'
$ws.Cells.Item(33, 3).Value = '13:00'
$ws.Cells.Item(33, 4).Value = 'Bathroom'
$ws.Cells.Item(33, 5).Value = '86.2%'
$ws.Cells.Item(33, 6).Value = 'Active'

$newRange.ClearFormats()

# --- Proximity sheet: append rows 41-43 ---
$ws = $wb.Worksheets.Item("Proximity")
$newRange = $ws.Range("A41:F43")
$newRange.NumberFormat = "@"

$ws.Cells.Item(41, 1).Value = '2026-02-01'
$ws.Cells.Item(41, 2).Value = '13:48:24'
$ws.Cells.Item(41, 3).Value = '13:00'
$ws.Cells.Item(41, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(41, 5).Value = 'EXIT'
$ws.Cells.Item(41, 6).Value = 'User EXITED Living Room Main Door'

$ws.Cells.Item(42, 1).Value = '2026-02-01'
$ws.Cells.Item(42, 2).Value = '13:48:40'
$ws.Cells.Item(42, 3).Value = '13:00'
$ws.Cells.Item(42, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(42, 5).Value = 'ENTER'
$ws.Cells.Item(42, 6).Value = 'User ENTERED Living Room Main Door'

$ws.Cells.Item(43, 1).Value = '2026-02-01'
$ws.Cells.Item(43, 2).Value = '13:49:11'
$ws.Cells.Item(43, 3).Value = '13:00'
$ws.Cells.Item(43, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(43, 5).Value = 'EXIT'
$ws.Cells.Item(43, 6).Value = 'User EXITED Living Room Main Door'

$newRange.ClearFormats()

# --- Camera sheet: append rows 19-20 ---
$ws = $wb.Worksheets.Item("Camera")
$newRange = $ws.Range("A19:F20")
$newRange.NumberFormat = "@"

$ws.Cells.Item(19, 1).Value = '2026-02-01'
$ws.Cells.Item(19, 2).Value = '13:48:24'
$ws.Cells.Item(19, 3).Value = '13:00'
$ws.Cells.Item(19, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(19, 5).Value = 'Image Received'
$ws.Cells.Item(19, 6).Value = 'Active'

$ws.Cells.Item(20, 1).Value = '2026-02-01'
$ws.Cells.Item(20, 2).Value = '13:49:11'
$ws.Cells.Item(20, 3).Value = '13:00'
$ws.Cells.Item(20, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(20, 5).Value = 'Image Received'
$ws.Cells.Item(20, 6).Value = 'Active'

$newRange.ClearFormats()
